$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
# Row 124
$ws.Cells.Item(124, 1).Value = 45910.4328125
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x01,0x7c"
$ws.Cells.Item(124, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(124, 4).Value = "0x01,0x08"
$ws.Cells.Item(124, 5).Value = "0x14"
$ws.Cells.Item(124, 6).Value = 380
$ws.Cells.Item(124, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(124, 8).Value = 264
$ws.Cells.Item(124, 9).Value = 14

# Row 125
$ws.Cells.Item(125, 1).Value = 45911.43440972222
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x01,0x7c"
$ws.Cells.Item(125, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(125, 4).Value = "0x01,0x08"
$ws.Cells.Item(125, 5).Value = "0x14"
$ws.Cells.Item(125, 6).Value = 380
$ws.Cells.Item(125, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(125, 8).Value = 264
$ws.Cells.Item(125, 9).Value = 14

# Row 126
$ws.Cells.Item(126, 1).Value = 45912.43680555555
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x01,0x7c"
$ws.Cells.Item(126, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(126, 4).Value = "0x01,0x08"
$ws.Cells.Item(126, 5).Value = "0x14"
$ws.Cells.Item(126, 6).Value = 380
$ws.Cells.Item(126, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item(126, 8).Value = 264
$ws.Cells.Item(126, 9).Value = 14

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
# Row 124
$ws.Cells.Item(124, 1).Value = 45910.4328125
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x01,0x7c"
$ws.Cells.Item(124, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(124, 4).Value = "0x01,0x10"
$ws.Cells.Item(124, 5).Value = "0xe"
$ws.Cells.Item(124, 6).Value = 380
$ws.Cells.Item(124, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(124, 8).Value = 272
$ws.Cells.Item(124, 9).Value = 14

# Row 125
$ws.Cells.Item(125, 1).Value = 45911.43440972222
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x01,0x7c"
$ws.Cells.Item(125, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(125, 4).Value = "0x01,0x0C"
$ws.Cells.Item(125, 5).Value = "0xe"
$ws.Cells.Item(125, 6).Value = 380
$ws.Cells.Item(125, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(125, 8).Value = 268
$ws.Cells.Item(125, 9).Value = 14

# Row 126
$ws.Cells.Item(126, 1).Value = 45912.43680555555
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x01,0x7c"
$ws.Cells.Item(126, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(126, 4).Value = "0x01,0x0C"
$ws.Cells.Item(126, 5).Value = "0xe"
$ws.Cells.Item(126, 6).Value = 380
$ws.Cells.Item(126, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(126, 8).Value = 268
$ws.Cells.Item(126, 9).Value = 14

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
# Row 124
$ws.Cells.Item(124, 1).Value = 45910.4328125
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x00,0x82"
$ws.Cells.Item(124, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(124, 4).Value = "0x00,0x70"
$ws.Cells.Item(124, 5).Value = "0x7"
$ws.Cells.Item(124, 6).Value = 130
$ws.Cells.Item(124, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(124, 8).Value = 112
$ws.Cells.Item(124, 9).Value = 7

# Row 125
$ws.Cells.Item(125, 1).Value = 45911.43440972222
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x00,0x82"
$ws.Cells.Item(125, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(125, 4).Value = "0x00,0x70"
$ws.Cells.Item(125, 5).Value = "0x7"
$ws.Cells.Item(125, 6).Value = 130
$ws.Cells.Item(125, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(125, 8).Value = 112
$ws.Cells.Item(125, 9).Value = 7

# Row 126
$ws.Cells.Item(126, 1).Value = 45912.43680555555
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x00,0x82"
$ws.Cells.Item(126, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(126, 4).Value = "0x00,0x70"
$ws.Cells.Item(126, 5).Value = "0x7"
$ws.Cells.Item(126, 6).Value = 130
$ws.Cells.Item(126, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item(126, 8).Value = 112
$ws.Cells.Item(126, 9).Value = 7

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
# Row 124
$ws.Cells.Item(124, 1).Value = 45910.4328125
$ws.Cells.Item(124, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(124, 2).Value = "0x00,0x82"
$ws.Cells.Item(124, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(124, 4).Value = "0x00,0x6E"
$ws.Cells.Item(124, 5).Value = "0x3"
$ws.Cells.Item(124, 6).Value = 130
$ws.Cells.Item(124, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(124, 8).Value = 110
$ws.Cells.Item(124, 9).Value = 3

# Row 125
$ws.Cells.Item(125, 1).Value = 45911.43440972222
$ws.Cells.Item(125, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 2).Value = "0x00,0x82"
$ws.Cells.Item(125, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(125, 4).Value = "0x00,0x6E"
$ws.Cells.Item(125, 5).Value = "0x3"
$ws.Cells.Item(125, 6).Value = 130
$ws.Cells.Item(125, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(125, 8).Value = 110
$ws.Cells.Item(125, 9).Value = 3

# Row 126
$ws.Cells.Item(126, 1).Value = 45912.43680555555
$ws.Cells.Item(126, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 2).Value = "0x00,0x82"
$ws.Cells.Item(126, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(126, 4).Value = "0x00,0x6D"
$ws.Cells.Item(126, 5).Value = "0x3"
$ws.Cells.Item(126, 6).Value = 130
$ws.Cells.Item(126, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item(126, 8).Value = 109
$ws.Cells.Item(126, 9).Value = 3
